# Update the gdp_calibrate sheet's "value" column (C2:C5) so the figures
# are more in line with the previous macro calculations, and leave the
# workbook with that sheet active/selected (as it was when the edit was
# made), matching the cell that was last touched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("gdp_calibrate")

# Make this the active sheet (updates workbook.xml activeTab + this
# sheet's tabSelected, and clears tabSelected on the previously active
# sheet).
$ws.Activate()

$ws.Range("C2").Value = 500
$ws.Range("C3").Value = 1000
$ws.Range("C4").Value = 2000
$ws.Range("C5").Value = 3000

# Leave the selection on the last edited cell.
$ws.Range("C5").Select()
